$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.9687
$ws.Range("A9").Value = -20.50169999999997
$ws.Range("E11").Value = 13.345
$ws.Range("A18").Value = -23.02440000000001
$ws.Range("A20").Value = -22.01860000000003
